$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N1").Value = "allowable_days_lag"
$ws.Range("O1").Value = "reviewed_at"
$ws.Range("P1").Value = "escalated_to"

$ws.Range("P3").Value = "Divisional perfomance meeting"
$ws.Range("O3").Value = "Service performance meeting"
